# Updates cryptocurrency price/volume figures (and restores the original
# row order for a couple of coin pairs that the source feed re-sorted) to
# match the latest scrape, per the commit:
#   "Updated cryptos list on Wed Jul 17 16:11:11 UTC 2024 with GitHub Actions"
#
# Values are written with a leading apostrophe so Excel stores them as text
# (matching the workbook's original inlineStr cells) instead of silently
# reinterpreting price strings such as "1.00" or "64.637.34" as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.637.34"
$ws.Range("E2").Value = "'  -0.24%  "

$ws.Range("D3").Value = "'3.428.61"
$ws.Range("E3").Value = "'  -1.39%  "

$ws.Range("E4").Value = "'  -0.08%  "

$ws.Range("D5").Value = "'571.27"
$ws.Range("E5").Value = "'  -0.82%  "

$ws.Range("D6").Value = "'157.59"
$ws.Range("E6").Value = "'  -1.78%  "

$ws.Range("E7").Value = "'  +0.11%  "

$ws.Range("D8").Value = "'0.599"
$ws.Range("E8").Value = "'  +2.34%  "

$ws.Range("D9").Value = "'3.434.43"
$ws.Range("E9").Value = "'  -1.65%  "

$ws.Range("D10").Value = "'7.24"
$ws.Range("E10").Value = "'  -2.08%  "

$ws.Range("D11").Value = "'0.122"
$ws.Range("E11").Value = "'  -2.71%  "

$ws.Range("D12").Value = "'0.445"
$ws.Range("E12").Value = "'  -0.10%  "

$ws.Range("D13").Value = "'4.028.90"
$ws.Range("E13").Value = "'  -1.38%  "

$ws.Range("E14").Value = "'  +0.09%  "

$ws.Range("D15").Value = "'0.0000189"
$ws.Range("E15").Value = "'  -3.72%  "

$ws.Range("D16").Value = "'27.90"
$ws.Range("E16").Value = "'  -1.83%  "

$ws.Range("D17").Value = "'64.747.22"
$ws.Range("E17").Value = "'  -0.11%  "

$ws.Range("D18").Value = "'3.447.91"
$ws.Range("E18").Value = "'  -1.38%  "

$ws.Range("D19").Value = "'6.35"
$ws.Range("E19").Value = "'  -1.46%  "

$ws.Range("D20").Value = "'14.09"
$ws.Range("E20").Value = "'  -2.25%  "

$ws.Range("D21").Value = "'377.13"
$ws.Range("E21").Value = "'  -4.12%  "

$ws.Range("D22").Value = "'8.05"
$ws.Range("E22").Value = "'  -2.97%  "

$ws.Range("D23").Value = "'0.552"
$ws.Range("E23").Value = "'  +1.15%  "

$ws.Range("B24").Value = "'Dai"
$ws.Range("C24").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "'  -0.44%  "

$ws.Range("B25").Value = "'Litecoin"
$ws.Range("C25").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'72.14"
$ws.Range("E25").Value = "'  -2.30%  "

$ws.Range("D26").Value = "'0.0000118"
$ws.Range("E26").Value = "'  -2.54%  "

$ws.Range("D27").Value = "'9.98"
$ws.Range("E27").Value = "'  +3.85%  "

$ws.Range("E28").Value = "'  -2.15%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "'  +0.05%  "

$ws.Range("D30").Value = "'1.48"
$ws.Range("E30").Value = "'  +1.52%  "

$ws.Range("D31").Value = "'6.09"
$ws.Range("E31").Value = "'  -1.48%  "

$ws.Range("D32").Value = "'2.03"
$ws.Range("E32").Value = "'  -0.33%  "

$ws.Range("D33").Value = "'23.33"
$ws.Range("E33").Value = "'  -1.67%  "

$ws.Range("D34").Value = "'7.19"
$ws.Range("E34").Value = "'  +2.21%  "

$ws.Range("E35").Value = "'  +6.66%  "

$ws.Range("D36").Value = "'159.77"
$ws.Range("E36").Value = "'  -0.84%  "

$ws.Range("D37").Value = "'1.89"
$ws.Range("E37").Value = "'  -0.09%  "

$ws.Range("D38").Value = "'0.0765"
$ws.Range("E38").Value = "'  -1.56%  "

$ws.Range("D39").Value = "'27.07"
$ws.Range("E39").Value = "'  -2.46%  "

$ws.Range("B40").Value = "'RenderToken"
$ws.Range("C40").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'6.79"
$ws.Range("E40").Value = "'  +2.02%  "

$ws.Range("B41").Value = "'Filecoin"
$ws.Range("C41").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'4.67"
$ws.Range("E41").Value = "'  +4.77%  "

$ws.Range("D42").Value = "'2.845.63"
$ws.Range("E42").Value = "'  -3.00%  "

$ws.Range("D43").Value = "'42.73"
$ws.Range("E43").Value = "'  -0.25%  "

$ws.Range("D44").Value = "'0.0315"
$ws.Range("E44").Value = "'  -1.29%  "

$ws.Range("D45").Value = "'0.774"
$ws.Range("E45").Value = "'  -0.55%  "

$ws.Range("D46").Value = "'25.79"
$ws.Range("E46").Value = "'  +8.41%  "

$ws.Range("B47").Value = "'ONDO"
$ws.Range("C47").Value = "'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "'1.08"
$ws.Range("E47").Value = "'  -2.95%  "

$ws.Range("B48").Value = "'Bittensor"
$ws.Range("C48").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "'316.26"
$ws.Range("E48").Value = "'  +6.24%  "

$ws.Range("E49").Value = "'  +0.98%  "

$ws.Range("D50").Value = "'0.863"
$ws.Range("E50").Value = "'  +0.49%  "

$ws.Range("D51").Value = "'6.56"
$ws.Range("E51").Value = "'  +0.33%  "
